$d = $word.ActiveDocument

# Paragraph 1 (FirstParagraph style): remove leading "**" and make bold
$d.Content.Find.Execute("**FOR IMMEDIATE RELEASE:", $true, $false, $false, $false, $false,
                         $true, 1, $false, "FOR IMMEDIATE RELEASE:", 2)

# Paragraph 3 (BodyText): remove trailing "**"
$d.Content.Find.Execute("Linux User Group Meeting**", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Linux User Group Meeting", 2)

# Now apply bold formatting to the three runs involved (text only, not the
# paragraph mark, so w:pPr/w:rPr stays untouched).
foreach ($p in $d.Paragraphs) {
    $pr = $p.Range
    $t = $pr.Text.TrimEnd("`r", "`n", [char]7)
    if ($t -eq "FOR IMMEDIATE RELEASE:" -or
        $t -eq "DESIRED PUBLICATION DATE: May 1, 2019" -or
        $t -eq "Linux User Group Meeting") {
        $r = $d.Range($pr.Start, $pr.End - 1)
        $r.Font.Bold = 1
    }
}
